$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain numeric-looking strings (e.g. "29.878.25",
# "0.0610"); Excel would otherwise auto-convert these to numbers on entry,
# silently dropping trailing zeros / the thousands-style dots. A leading
# apostrophe forces literal text entry (Excel's normal "text cell" syntax)
# while leaving the cell format as "General", exactly like the source.
$ws.Range("D2").Value = "'29.878.25"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "'1.624.83"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'214.63"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "'0.521"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").Value = "'29.66"
$ws.Range("E8").Value = "  +9.29%  "
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("D10").Value = "'0.0610"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "'1.858.24"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "'1.621.52"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +6.14%  "
$ws.Range("D15").Value = "'3.90"
$ws.Range("E15").Value = "  +4.69%  "
$ws.Range("D16").Value = "'29.943.81"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "'8.87"
$ws.Range("E17").Value = "  +16.66%  "
$ws.Range("D18").Value = "'64.62"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").Value = "'243.51"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "'0.0₃0705"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  +3.41%  "
$ws.Range("D23").Value = "'9.61"
$ws.Range("E23").Value = "  +4.02%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "'157.46"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").Value = "'15.65"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("D28").Value = "'6.60"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'0.0489"
$ws.Range("E30").Value = "  +2.93%  "
$ws.Range("E31").Value = "  +5.21%  "
$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("D33").Value = "'3.22"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").Value = "'1.420.68"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +6.52%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'2.87"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("D40").Value = "'0.558"
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("E43").Value = "  +3.93%  "
$ws.Range("D44").Value = "'54.62"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'69.22"
$ws.Range("E45").Value = "  +4.93%  "
$ws.Range("E46").Value = "  +13.55%  "
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("D49").Value = "'1.766.58"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").Value = "'88.73"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'0.0₆0109"
$ws.Range("E51").Value = "  +3.54%  "
